$wb = $excel.ActiveWorkbook

# Rename shared strings: "light goods" -> "van", "heavy goods" -> "lorry"
# These labels appear in column A (rows 3 and 4) of every worksheet in the workbook.
foreach ($ws in $wb.Worksheets) {
    $ws.Range("A3").Value = "van"
    $ws.Range("A4").Value = "lorry"
}

# Sheet: mean
$ws = $wb.Worksheets.Item("mean")
$ws.Cells.Item(2, 2).Value = 26.202543800299
$ws.Cells.Item(2, 3).Value = 44.09650293919952
$ws.Cells.Item(2, 4).Value = 33.75668622226373
$ws.Cells.Item(2, 5).Value = 42.93465177040636
$ws.Cells.Item(2, 6).Value = 36.96310806496159
$ws.Cells.Item(3, 2).Value = 11.979800505986175
$ws.Cells.Item(3, 3).Value = 20.08626932259857
$ws.Cells.Item(3, 4).Value = 14.616786362421816
$ws.Cells.Item(3, 5).Value = 15.223735748636328
$ws.Cells.Item(3, 6).Value = 14.998258034415239
$ws.Cells.Item(4, 2).Value = 15.856086415537916
$ws.Cells.Item(4, 3).Value = 23.05559686613122
$ws.Cells.Item(4, 4).Value = 45.0463723682924
$ws.Cells.Item(4, 5).Value = 31.591169667806824
$ws.Cells.Item(4, 6).Value = 21.603393504324313
$ws.Cells.Item(5, 2).Value = 98.65683464269108
$ws.Cells.Item(5, 3).Value = 145.429751456658
$ws.Cells.Item(5, 4).Value = 66.70418367992117
$ws.Cells.Item(5, 5).Value = 40.176610920065656
$ws.Cells.Item(5, 6).Value = 71.55295026329189
$ws.Cells.Item(6, 2).Value = 65.12139240093677
$ws.Cells.Item(6, 3).Value = 154.9654757533907
$ws.Cells.Item(6, 4).Value = 97.60917776436251
$ws.Cells.Item(6, 5).Value = 113.26054009644454
$ws.Cells.Item(6, 6).Value = 114.87668805639032
$ws.Cells.Item(7, 2).Value = 131.30229311345965
$ws.Cells.Item(7, 3).Value = 39.14559352429843
$ws.Cells.Item(7, 4).Value = 22.080318155527102
$ws.Cells.Item(7, 5).Value = 11.761666513867201
$ws.Cells.Item(7, 6).Value = 19.276966978164797

# Sheet: median
$ws = $wb.Worksheets.Item("median")
$ws.Cells.Item(2, 2).Value = 26.194920094252602
$ws.Cells.Item(2, 3).Value = 44.091468864582026
$ws.Cells.Item(2, 4).Value = 33.753736957622
$ws.Cells.Item(2, 5).Value = 42.92540443460415
$ws.Cells.Item(2, 6).Value = 36.95874826403663
$ws.Cells.Item(3, 2).Value = 11.979968662931022
$ws.Cells.Item(3, 3).Value = 20.035118004598893
$ws.Cells.Item(3, 4).Value = 14.593963331501623
$ws.Cells.Item(3, 5).Value = 15.189440058753455
$ws.Cells.Item(3, 6).Value = 14.969016943941593
$ws.Cells.Item(4, 2).Value = 15.724771196106463
$ws.Cells.Item(4, 3).Value = 22.878280410724123
$ws.Cells.Item(4, 4).Value = 44.611443451096996
$ws.Cells.Item(4, 5).Value = 31.160846752967945
$ws.Cells.Item(4, 6).Value = 21.565950139674555
$ws.Cells.Item(5, 2).Value = 98.21206094992274
$ws.Cells.Item(5, 3).Value = 145.08599046848497
$ws.Cells.Item(5, 4).Value = 66.28881157051043
$ws.Cells.Item(5, 5).Value = 40.03890310878499
$ws.Cells.Item(5, 6).Value = 71.2206369946636
$ws.Cells.Item(6, 2).Value = 64.95334149931585
$ws.Cells.Item(6, 3).Value = 155.03647491223097
$ws.Cells.Item(6, 4).Value = 97.7172285751854
$ws.Cells.Item(6, 5).Value = 112.76851186129974
$ws.Cells.Item(6, 6).Value = 114.71660410392971
$ws.Cells.Item(7, 2).Value = 129.45459592633054
$ws.Cells.Item(7, 3).Value = 39.000853436284544
$ws.Cells.Item(7, 4).Value = 22.04388813465009
$ws.Cells.Item(7, 5).Value = 11.718989402550385
$ws.Cells.Item(7, 6).Value = 19.324829606391432

# Sheet: lower 5
$ws = $wb.Worksheets.Item("lower 5")
$ws.Cells.Item(2, 2).Value = 25.76628260602799
$ws.Cells.Item(2, 3).Value = 43.436185552459655
$ws.Cells.Item(2, 4).Value = 33.13161648377678
$ws.Cells.Item(2, 5).Value = 42.4017319360135
$ws.Cells.Item(2, 6).Value = 36.69711379300097
$ws.Cells.Item(3, 2).Value = 10.826752801566464
$ws.Cells.Item(3, 3).Value = 18.154983076599645
$ws.Cells.Item(3, 4).Value = 13.117471289749375
$ws.Cells.Item(3, 5).Value = 13.868256260358185
$ws.Cells.Item(3, 6).Value = 13.820720180098029
$ws.Cells.Item(4, 2).Value = 13.377067949152837
$ws.Cells.Item(4, 3).Value = 18.971510615306233
$ws.Cells.Item(4, 4).Value = 36.74567878931382
$ws.Cells.Item(4, 5).Value = 25.415479459307747
$ws.Cells.Item(4, 6).Value = 18.556388231310656
$ws.Cells.Item(5, 2).Value = 82.68497931828918
$ws.Cells.Item(5, 3).Value = 125.92408911487254
$ws.Cells.Item(5, 4).Value = 55.36731348536185
$ws.Cells.Item(5, 5).Value = 34.83687979249456
$ws.Cells.Item(5, 6).Value = 62.87367217408012
$ws.Cells.Item(6, 2).Value = 56.16966347347308
$ws.Cells.Item(6, 3).Value = 139.55736053364095
$ws.Cells.Item(6, 4).Value = 83.64724564043821
$ws.Cells.Item(6, 5).Value = 102.95688728545731
$ws.Cells.Item(6, 6).Value = 104.99767060416099
$ws.Cells.Item(7, 2).Value = 102.73731329767088
$ws.Cells.Item(7, 3).Value = 33.33918479890467
$ws.Cells.Item(7, 4).Value = 18.23784077763773
$ws.Cells.Item(7, 5).Value = 10.23786295290008
$ws.Cells.Item(7, 6).Value = 16.942836914282125

# Sheet: upper 95
$ws = $wb.Worksheets.Item("upper 95")
$ws.Cells.Item(2, 2).Value = 26.676257848231277
$ws.Cells.Item(2, 3).Value = 44.791012337355475
$ws.Cells.Item(2, 4).Value = 34.41121883716436
$ws.Cells.Item(2, 5).Value = 43.452729445079726
$ws.Cells.Item(2, 6).Value = 37.24292258175648
$ws.Cells.Item(3, 2).Value = 13.099440626927073
$ws.Cells.Item(3, 3).Value = 22.016400515194125
$ws.Cells.Item(3, 4).Value = 16.21977271522888
$ws.Cells.Item(3, 5).Value = 16.709659252422668
$ws.Cells.Item(3, 6).Value = 16.168087848234407
$ws.Cells.Item(4, 2).Value = 18.623478757543992
$ws.Cells.Item(4, 3).Value = 27.24340296113507
$ws.Cells.Item(4, 4).Value = 54.478459034106535
$ws.Cells.Item(4, 5).Value = 39.046507545587616
$ws.Cells.Item(4, 6).Value = 24.863672134453743
$ws.Cells.Item(5, 2).Value = 115.89526273505956
$ws.Cells.Item(5, 3).Value = 166.21187842975337
$ws.Cells.Item(5, 4).Value = 79.55000104943484
$ws.Cells.Item(5, 5).Value = 46.20857579848022
$ws.Cells.Item(5, 6).Value = 80.93058648539562
$ws.Cells.Item(6, 2).Value = 73.98397550733814
$ws.Cells.Item(6, 3).Value = 170.7738639608777
$ws.Cells.Item(6, 4).Value = 112.71288391904096
$ws.Cells.Item(6, 5).Value = 125.13836024649409
$ws.Cells.Item(6, 6).Value = 125.01420327846127
$ws.Cells.Item(7, 2).Value = 166.3117564882687
$ws.Cells.Item(7, 3).Value = 45.391583801194436
$ws.Cells.Item(7, 4).Value = 26.198514495852507
$ws.Cells.Item(7, 5).Value = 13.334833665979955
$ws.Cells.Item(7, 6).Value = 21.554148245906003
